# Add two new rows (10 and 11) documenting the XGBoosting income-prediction
# models, mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (A=plain vertical-centered, B=wrap vertical-centered, C=plain
# vertical-centered) is the closest existing template for the new rows, so
# copy its formatting into A10:C10 and A11:C11 before writing values.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("A9:C9").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(10, 1).Value = "Predicción Ingreso"
$ws.Cells.Item(10, 2).Value = @'
Ingtotug~Dominio + Depto + P5010 + 
                      N_cuartos_hog + Nper + nmenores_5 + nmenores_6_11 + 
                      nmenores_12_17 + nocupados + nincapacitados + ntrabajo_menores + 
                      Head_Mujer + Head_Afiliado_SS + P5140 + Npersug +
                      Head_exper_ult_trab + Head_Rec_alimento + Head_Rec_subsidio + 
                      Head_Rec_vivienda + Head_Ocupacion + maxEducLevel + Head_Primas +
                      Head_Segundo_trabajo + DormitorXpersona + Ln_Cuota + Head_Oficio +
                      Ln_Pago_arrien + nmujeres + Ocup_vivienda + 
                      Head_Cot_pension + Cabecera
'@
$ws.Cells.Item(10, 3).Value = "XGBoosting"

$ws.Cells.Item(11, 1).Value = "Predicción Ingreso"
$ws.Cells.Item(11, 2).Value = @'
Ln_Ing_tot_hogar~Dominio + Depto + P5010 + 
                        N_cuartos_hog + Nper + nmenores_5 + nmenores_6_11 + 
                        nmenores_12_17 + nocupados + nincapacitados + ntrabajo_menores + 
                        Head_Mujer + Head_Afiliado_SS + P5140 + Npersug +
                        Head_exper_ult_trab + Head_Rec_alimento + Head_Rec_subsidio + 
                        Head_Rec_vivienda + Head_Ocupacion + maxEducLevel + Head_Primas +
                        Head_Segundo_trabajo + DormitorXpersona + Ln_Cuota + Head_Oficio +
                        Ln_Pago_arrien + nmujeres + Ocup_vivienda + 
                        Head_Cot_pension + Cabecera
'@
$ws.Cells.Item(11, 3).Value = "XGBoosting"

# The long, multi-line formula text needs a taller row - match the height
# computed by Excel for this amount of wrapped text.
$ws.Rows.Item(10).RowHeight = 144
$ws.Rows.Item(11).RowHeight = 144

$ws.Range("C12").Select()
